$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("County")

function Set-TextValue($rangeAddr, $text) {
    $cell = $ws.Range($rangeAddr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue "B30" '0.00%'
Set-TextValue "C30" '$0'
Set-TextValue "D30" '0.00%'
Set-TextValue "E30" '0.00%'
Set-TextValue "F30" '0.00%'

Set-TextValue "B31" '0.00%'
Set-TextValue "C31" '$0'
Set-TextValue "D31" '0.00%'
Set-TextValue "E31" '0.00%'
Set-TextValue "F31" '0.00%'
